$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Free Guy"
$ws.Range("B2").Value = "A bank teller called Guy realizes he is a background character in an open world video game called Free City that will soon go offline."

$ws.Range("A3").Value = "Eternals "
$ws.Range("B3").Value = "The Eternals are a team of ancient aliens who have been living on Earth in secret for thousands of years. When an unexpected tragedy forces them out of the shadows, they are forced to reunite against mankind’s most ancient enemy, the Deviants."

$ws.Range("A4").Value = "Spider-Man: No Way Home "
$ws.Range("B4").Value = "Peter Parker is unmasked and no longer able to separate his normal life from the high-stakes of being a Super Hero. When he asks for help from Doctor Strange the stakes become even more dangerous, forcing him to discover what it truly means to be Spider-Man."
